$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "testing"
$ws.Range("B4").Value = "test"
